$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33 (G33=5512)
$ws.Range("H33").Value = 221.53334
$ws.Range("I33").Value = 231.5
$ws.Range("J33").Value = 82
$ws.Range("K33").Value = 231.5
$ws.Range("L33").Value = 82
$ws.Range("M33").Value = -2.5
$ws.Range("N33").Value = -540
# Row 45 (G45=4585)
$ws.Range("H45").Value = 6633.4
$ws.Range("J45").Value = 2800
$ws.Range("L45").Value = 8400
$ws.Range("N45").Value = -8784
# Row 62 (G62=27781)
$ws.Range("H62").Value = 3189.7778
$ws.Range("I62").Value = 3234.5
$ws.Range("J62").Value = 3100.3333
$ws.Range("K62").Value = 3234.5
$ws.Range("L62").Value = 3100.3333
$ws.Range("M62").Value = -2610.5
$ws.Range("N62").Value = -4348.3333
# Row 65 (G65=27781)
$ws.Range("H65").Value = 3189.7778
$ws.Range("I65").Value = 3234.5
$ws.Range("J65").Value = 3100.3333
$ws.Range("K65").Value = 16172.5
$ws.Range("L65").Value = 15501.6665
$ws.Range("M65").Value = -13052.5
$ws.Range("N65").Value = -21741.6665
# Row 129 (G129=36115)
$ws.Range("H129").Value = 1153.7561
$ws.Range("J129").Value = 1166.7693
$ws.Range("L129").Value = 3500.3079
$ws.Range("N129").Value = -13500.3079
# Row 132 (G132=44049)
$ws.Range("H132").Value = 1619.2407
$ws.Range("I132").Value = 1619.2407
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4857.7221
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2327.7221
$ws.Range("N132").ClearContents()
# Row 134 (G134=41997)
$ws.Range("H134").Value = 116702.71
$ws.Range("J134").Value = 116702.71
$ws.Range("L134").Value = 116702.71
$ws.Range("N134").Value = -126842.71
# Row 136 (G136=42164)
$ws.Range("H136").Value = 62797.145
$ws.Range("J136").Value = 62797.145
$ws.Range("L136").Value = 62797.145
$ws.Range("N136").Value = -72997.14499999999
# Row 137 (G137=44013)
$ws.Range("H137").Value = 2716.7917
$ws.Range("I137").Value = 2512.4375
$ws.Range("J137").Value = 3125.5
$ws.Range("K137").Value = 7537.3125
$ws.Range("L137").Value = 9376.5
$ws.Range("M137").Value = -4987.3125
$ws.Range("N137").Value = -14476.5
# Row 138 (G138=44169)
$ws.Range("H138").Value = 6459714.5
$ws.Range("I138").Value = 12504541
$ws.Range("J138").Value = 11899
$ws.Range("K138").Value = 37513623
$ws.Range("L138").Value = 35697
$ws.Range("M138").Value = -37508483
$ws.Range("N138").Value = -45977

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32 (G32=44147)
$ws.Range("H32").Value = 51391.434
$ws.Range("I32").Value = 40912.69
$ws.Range("K32").Value = 40912.69
$ws.Range("M32").Value = -40625.69
# Row 61 (G61=43999)
$ws.Range("H61").Value = 2326.0908
$ws.Range("I61").Value = 1726.2693
$ws.Range("J61").Value = 4554
$ws.Range("K61").Value = 1726.2693
$ws.Range("L61").Value = 4554
$ws.Range("M61").Value = -1514.2693
$ws.Range("N61").Value = -4978
# Row 74 (G74=44000)
$ws.Range("H74").Value = 1285.1613
$ws.Range("I74").Value = 1244.6666
$ws.Range("J74").Value = 2500
$ws.Range("K74").Value = 1244.6666
$ws.Range("L74").Value = 2500
$ws.Range("M74").Value = -370.6666
$ws.Range("N74").Value = -4248
# Row 77 (G77=44000)
$ws.Range("H77").Value = 1285.1613
$ws.Range("I77").Value = 1244.6666
$ws.Range("J77").Value = 2500
$ws.Range("K77").Value = 6223.333000000001
$ws.Range("L77").Value = 12500
$ws.Range("M77").Value = -1855.333000000001
$ws.Range("N77").Value = -21236
# Row 119 (G119=26287)
$ws.Range("H119").Value = 38799.5
$ws.Range("J119").Value = 38799.5
$ws.Range("L119").Value = 38799.5
$ws.Range("N119").Value = -48475.5
# Row 132 (G132=43997)
$ws.Range("H132").Value = 1776.6558
$ws.Range("I132").Value = 1543.6666
$ws.Range("J132").Value = 2964.9
$ws.Range("K132").Value = 4630.9998
$ws.Range("L132").Value = 8894.700000000001
$ws.Range("M132").Value = -2100.9998
$ws.Range("N132").Value = -13954.7
# Row 136 (G136=43999)
$ws.Range("H136").Value = 2326.0908
$ws.Range("I136").Value = 1726.2693
$ws.Range("J136").Value = 4554
$ws.Range("K136").Value = 5178.8079
$ws.Range("L136").Value = 13662
$ws.Range("M136").Value = -2628.8079
$ws.Range("N136").Value = -18762

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 109 (G109=27096)
$ws.Range("H109").Value = 25596.111
$ws.Range("J109").Value = 25596.111
$ws.Range("L109").Value = 25596.111
$ws.Range("N109").Value = -28370.111

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31 (G31=44023)
$ws.Range("H31").Value = 5252.7856
$ws.Range("I31").Value = 4452.5
$ws.Range("J31").Value = 7253.5
$ws.Range("K31").Value = 4452.5
$ws.Range("L31").Value = 7253.5
$ws.Range("M31").Value = -4157.5
$ws.Range("N31").Value = -7843.5
# Row 34 (G34=44023)
$ws.Range("H34").Value = 5252.7856
$ws.Range("I34").Value = 4452.5
$ws.Range("J34").Value = 7253.5
$ws.Range("K34").Value = 4452.5
$ws.Range("L34").Value = 7253.5
$ws.Range("M34").Value = -4250.5
$ws.Range("N34").Value = -7657.5

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5 (G5=43974)
$ws.Range("H5").Value = 1620.3334
$ws.Range("I5").Value = 1332.9474
$ws.Range("J5").Value = 1893.35
$ws.Range("K5").Value = 3998.8422
$ws.Range("L5").Value = 5680.049999999999
$ws.Range("M5").Value = -3886.8422
$ws.Range("N5").Value = -5904.049999999999
# Row 81 (G81=12843)
$ws.Range("H81").Value = 3940
$ws.Range("I81").Value = 1300
$ws.Range("J81").Value = 4600
$ws.Range("K81").Value = 3900
$ws.Range("L81").Value = 13800
$ws.Range("M81").Value = -2777
$ws.Range("N81").Value = -16046
# Row 84 (G84=12843)
$ws.Range("H84").Value = 3940
$ws.Range("I84").Value = 1300
$ws.Range("J84").Value = 4600
$ws.Range("K84").Value = 11700
$ws.Range("L84").Value = 41400
$ws.Range("M84").Value = -6084
$ws.Range("N84").Value = -52632
# Row 107 (G107=27838)
$ws.Range("H107").Value = 512.6
$ws.Range("J107").Value = 349.8
$ws.Range("L107").Value = 1049.4
$ws.Range("N107").Value = -4889.4
# Row 112 (G112=27855)
$ws.Range("H112").Value = 3794.8
$ws.Range("I112").Value = 1199.7778
$ws.Range("J112").Value = 4906.952
$ws.Range("K112").Value = 3599.3334
$ws.Range("L112").Value = 14720.856
$ws.Range("M112").Value = -2491.3334
$ws.Range("N112").Value = -16936.856
# Row 113 (G113=27843)
$ws.Range("H113").Value = 9910.416999999999
$ws.Range("I113").Value = 785.2
$ws.Range("J113").Value = 16428.428
$ws.Range("K113").Value = 2355.6
$ws.Range("L113").Value = 49285.284
$ws.Range("M113").Value = -185.6000000000004
$ws.Range("N113").Value = -53625.284
# Row 122 (G122=36078)
$ws.Range("H122").Value = 50010900
$ws.Range("J122").Value = 100021096
$ws.Range("L122").Value = 900189864
$ws.Range("N122").Value = -900194764
# Row 123 (G123=36037)
$ws.Range("H123").Value = 10250
$ws.Range("J123").Value = 13333.333
$ws.Range("L123").Value = 39999.999
$ws.Range("N123").Value = -44899.999
# Row 135 (G135=43974)
$ws.Range("H135").Value = 1620.3334
$ws.Range("I135").Value = 1332.9474
$ws.Range("J135").Value = 1893.35
$ws.Range("K135").Value = 11996.5266
$ws.Range("L135").Value = 17040.15
$ws.Range("M135").Value = -9461.526600000001
$ws.Range("N135").Value = -22110.15

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 119 (G119=26288)
$ws.Range("H119").Value = 39000
$ws.Range("J119").Value = 39000
$ws.Range("L119").Value = 39000
$ws.Range("N119").Value = -48676
# Row 133 (G133=41903)
$ws.Range("H133").Value = 69748.11
$ws.Range("J133").Value = 69748.11
$ws.Range("L133").Value = 69748.11
$ws.Range("N133").Value = -74808.11
# Row 136 (G136=44060)
$ws.Range("H136").Value = 7058.8423
$ws.Range("I136").Value = 7807.533
$ws.Range("J136").Value = 4251.25
$ws.Range("K136").Value = 23422.599
$ws.Range("L136").Value = 12753.75
$ws.Range("M136").Value = -20872.599
$ws.Range("N136").Value = -17853.75

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 119 (G119=26289)
$ws.Range("H119").Value = 40000
$ws.Range("J119").Value = 40000
$ws.Range("L119").Value = 40000
$ws.Range("N119").Value = -49676
# Row 136 (G136=44031)
$ws.Range("H136").Value = 1816.45
$ws.Range("I136").Value = 1803.6451
$ws.Range("J136").Value = 1860.5555
$ws.Range("K136").Value = 5410.9353
$ws.Range("L136").Value = 5581.666499999999
$ws.Range("M136").Value = -2860.9353
$ws.Range("N136").Value = -10681.6665
